# Update Code to website
# Include website for missing packages.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top. The former row 1 ("Requires IT (Mike) to
# install", highlighted) shifts down to row 2, keeping its formatting.
$ws.Rows.Item(1).Insert()
$ws.Range("A1").Value = "Workbook created using Rpackage_status.R. Code needs to be rerun with each new install."

# Insert a new column D (to the right of the "Missing Package notes"
# column) for the package download-site links. Existing formatting on
# row 4 (header) / row 5 / row 10 (highlighted rows) is inherited
# automatically for the new column cells.
$ws.Columns.Item(4).Insert()
$ws.Columns.Item(4).ColumnWidth = 71.71

# Header row
$ws.Range("D4").Value = "Missing Package Download Site"

# Package download links
$ws.Range("D5").Value  = "https://cran.r-project.org/web/packages/processx/index.html"
$ws.Range("D6").Value  = "https://cran.r-project.org/web/packages/callr/index.html"
$ws.Range("D7").Value  = "https://cran.r-project.org/web/packages/reprex/index.html"
$ws.Range("D8").Value  = "https://cran.r-project.org/web/packages/rvest/index.html"
$ws.Range("D9").Value  = "https://cran.r-project.org/web/packages/tidyverse/index.html"
$ws.Range("D10").Value = "https://cran.r-project.org/web/packages/rgdal/index.html"
$ws.Range("D11").Value = "https://cran.r-project.org/web/packages/tigris/index.html"
$ws.Range("D12").Value = "https://cran.r-project.org/web/packages/choroplethr/index.html"
